$p = $ppt.ActivePresentation

# --- Slide 1 : subtitle date run "April XX, 2023" -> "April 12, 2023" ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$dateRun1 = $sh1.TextFrame.TextRange.Characters(1, 14)
$dateRun1.Text = "April 12, 2023"

# --- Slide 2 : title "Problem Statement: We need data!" -> "...We need legal data!" ---
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Problem Statement: We need legal data!"

# --- Slide 2 : body paragraph gains a second question ---
$bodyShape = $s2.Shapes.Item(2)
$bodyShape.TextFrame.TextRange.Text = "We want to train an AI model to extract information from manually-scanned PDFs. We only have so many documents to train from, and they don’t represent the full range of possible smudges and corruptions. What can we do to make our model more robust? What can we do to shield our data-providing customers from the model’s training?"

# --- Slide 6 : title's last run " ethical? It's fake, but…" -> " ethical? It's faked, but…" ---
$s6 = $p.Slides.Item(6)
$titleShape6 = $s6.Shapes.Item(1)
$fullTitleText = $titleShape6.TextFrame.TextRange.Text
$lastRun = $titleShape6.TextFrame.TextRange.Characters(25, $fullTitleText.Length - 24)
$lastRun.Text = " ethical? It’s faked, but…"

# --- Slide 6 : decorative bar shape reposition/resize ---
$barShape = $s6.Shapes.Item(5)
$barShape.Left = 0.002
$barShape.Width = 720.0

# --- Slide 7 : subtitle date run "April XX, 2023" -> "April 12, 2023" ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$dateRun7 = $sh7.TextFrame.TextRange.Characters(1, 14)
$dateRun7.Text = "April 12, 2023"
